$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): shift headers one column right and rename ---
$ws.Range("A1").Value = "Opcion"
$ws.Range("B1").Value = "Side"
$ws.Range("C1").Value = "Strike"
$ws.Range("D1").Value = "Price"
$ws.Range("E1").Value = "Cant"

# Clear the now-unused F1 cell (old 6th header column no longer exists) before
# touching the clipboard, since Clear() drops any pending Copy/Cut selection.
$ws.Range("F1").Clear() | Out-Null

# --- Data rows ---
$ws.Range("A2").Value = "GFGC100OCT"
$ws.Range("B2").Value = "C"
$ws.Range("C2").Value = 100
$ws.Range("D2").Value = 5
$ws.Range("E2").Value = 10

$ws.Range("A3").Value = "GFGC120OCT"
$ws.Range("B3").Value = "C"
$ws.Range("C3").Value = 120
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 3

$ws.Range("A4").Value = "GFGV80OCT"
$ws.Range("B4").Value = "V"
$ws.Range("C4").Value = 80
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = -5

# Apply the bold/centered/bordered header style (currently on B1) to the full
# header row plus column A of the data rows. PasteSpecial(Formats) reuses the
# existing cellXfs entry instead of minting a new (duplicate) style.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A1:E1").PasteSpecial(-4122) | Out-Null
$ws.Range("A2:A4").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false
